{"js": "const pairs = [\n  [\"2024-09-26 Thursday\", \"2024-09-27 Friday\"],\n  [\"92-83=\", \"60+12=\"],\n  [\"34+48=\", \"37+23=\"],\n  [\"86-76=\", \"68+2=\"],\n  [\"69-12=\", \"54+35=\"],\n  [\"93-63=\", \"62-51=\"],\n  [\"76-39=\", \"15+29=\"],\n  [\"56-31=\", \"85-19=\"],\n  [\"52+17=\", \"63-40=\"],\n  [\"99-85=\", \"48+15=\"],\n  [\"30+2=\", \"53+45=\"],\n  [\"65-33=\", \"58+27=\"],\n  [\"34+60=\", \"55+19=\"],\n  [\"97-93=\", \"7+55=\"],\n  [\"26+45=\", \"14+76=\"],\n  [\"70+5=\", \"65+16=\"],\n  [\"93-16=\", \"5+10=\"],\n  [\"6+13=\", \"24+59=\"],\n  [\"58-57=\", \"29-9=\"],\n  [\"66-8=\", \"91-86=\"],\n  [\"53+2=\", \"4+24=\"],\n  [\"66-47=\", \"99-2=\"],\n  [\"25+49=\", \"72-10=\"],\n  [\"19+35=\", \"25+37=\"],\n  [\"58-2=\", \"42-40=\"],\n  [\"44-10=\", \"80-79=\"],\n  [\"27+16=\", \"46+5=\"],\n  [\"26-7=\", \"83-16=\"],\n  [\"57+22=\", \"33+63=\"],\n  [\"61+12=\", \"3+63=\"],\n  [\"25+55=\", \"85-3=\"],\n  [\"91-72=\", \"65+11=\"],\n  [\"69-56=\", \"34-28=\"],\n  [\"73+18=\", \"7+63=\"],\n  [\"77-23=\", \"46-44=\"],\n  [\"16+82=\", \"1+75=\"],\n  [\"55-28=\", \"48-15=\"],\n  [\"95-63=\", \"96-49=\"],\n  [\"11-3=\", \"73+24=\"],\n  [\"39+54=\", \"61-53=\"],\n  [\"22+11=\", \"82-11=\"],\n  [\"23+11=\", \"36+17=\"],\n  [\"75+16=\", \"48-28=\"],\n  [\"20+70=\", \"54+18=\"],\n  [\"37+33=\", \"82-41=\"],\n  [\"74+9=\", \"0+46=\"],\n  [\"1+88=\", \"97-61=\"],\n  [\"59+23=\", \"79+6=\"],\n  [\"4+56=\", \"52-47=\"],\n  [\"32-12=\", \"45+32=\"],\n  [\"75-71=\", \"65+28=\"],\n  [\"83-76=\", \"68-11=\"],\n  [\"40+42=\", \"56-12=\"],\n  [\"96-4=\", \"2+77=\"],\n  [\"90-71=\", \"89-34=\"],\n  [\"94+2=\", \"46+0=\"],\n  [\"1+17=\", \"12+30=\"],\n  [\"24+23=\", \"51-40=\"],\n  [\"5+7=\", \"95-48=\"],\n  [\"82-20=\", \"25+29=\"],\n  [\"60-54=\", \"42-41=\"],\n  [\"23-16=\", \"37+19=\"],\n  [\"64-38=\", \"24+24=\"],\n  [\"9+15=\", \"87-20=\"],\n  [\"92-75=\", \"83-53=\"],\n  [\"79+19=\", \"91-54=\"],\n  [\"67-56=\", \"85-59=\"],\n  [\"2+66=\", \"36+52=\"],\n  [\"62-3=\", \"88-65=\"],\n  [\"30+49=\", \"78-42=\"],\n  [\"75-61=\", \"59-31=\"],\n  [\"27+31=\", \"12-11=\"],\n  [\"56-49=\", \"89-11=\"],\n  [\"48+6=\", \"74-8=\"],\n  [\"46-7=\", \"14+57=\"],\n  [\"52+14=\", \"61-18=\"],\n  [\"64-57=\", \"63+3=\"],\n  [\"28+32=\", \"53+18=\"],\n  [\"70+19=\", \"4+39=\"],\n  [\"4+93=\", \"68-42=\"],\n  [\"79-2=\", \"36+5=\"],\n  [\"32+67=\", \"49+9=\"],\n  [\"8+58=\", \"37+22=\"],\n  [\"1+87=\", \"49-48=\"],\n  [\"58-30=\", \"43+35=\"],\n  [\"76-5=\", \"39+55=\"],\n  [\"11+82=\", \"5+51=\"],\n  [\"78-52=\", \"49-33=\"],\n  [\"80-73=\", \"10+59=\"],\n  [\"46+49=\", \"16+45=\"],\n  [\"52+13=\", \"11+80=\"],\n  [\"49-44=\", \"79-64=\"],\n  [\"33-6=\", \"73+12=\"],\n  [\"20+15=\", \"24-17=\"],\n  [\"49+38=\", \"39-32=\"],\n  [\"58-29=\", \"70-13=\"],\n  [\"19-13=\", \"73-69=\"],\n  [\"66+3=\", \"0+54=\"],\n  [\"35+14=\", \"46+44=\"],\n  [\"41+40=\", \"45-16=\"],\n  [\"97-41=\", \"81-80=\"],\n];\n\nconst body = context.document.body;\n\n// Kick off a search for every old->new pair up front (batched), then\n// sync once so all results come back together.\nconst searchResults = pairs.map(([oldText]) =>\n  body.search(oldText, { matchCase: true, matchWholeWord: false })\n);\nsearchResults.forEach((r) => r.load(\"items,text\"));\nawait context.sync();\n\n// Replace each unique match with its new text. Every source string in\n// this worksheet is unique, so each search should resolve to exactly\n// one range.\nfor (let i = 0; i < pairs.length; i++) {\n  const [oldText, newText] = pairs[i];\n  const items = searchResults[i].items;\n  if (items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n  for (const item of items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$pairs = @(\n  @(\"2024-09-26 Thursday\", \"2024-09-27 Friday\"),\n  @(\"92-83=\", \"60+12=\"),\n  @(\"34+48=\", \"37+23=\"),\n  @(\"86-76=\", \"68+2=\"),\n  @(\"69-12=\", \"54+35=\"),\n  @(\"93-63=\", \"62-51=\"),\n  @(\"76-39=\", \"15+29=\"),\n  @(\"56-31=\", \"85-19=\"),\n  @(\"52+17=\", \"63-40=\"),\n  @(\"99-85=\", \"48+15=\"),\n  @(\"30+2=\", \"53+45=\"),\n  @(\"65-33=\", \"58+27=\"),\n  @(\"34+60=\", \"55+19=\"),\n  @(\"97-93=\", \"7+55=\"),\n  @(\"26+45=\", \"14+76=\"),\n  @(\"70+5=\", \"65+16=\"),\n  @(\"93-16=\", \"5+10=\"),\n  @(\"6+13=\", \"24+59=\"),\n  @(\"58-57=\", \"29-9=\"),\n  @(\"66-8=\", \"91-86=\"),\n  @(\"53+2=\", \"4+24=\"),\n  @(\"66-47=\", \"99-2=\"),\n  @(\"25+49=\", \"72-10=\"),\n  @(\"19+35=\", \"25+37=\"),\n  @(\"58-2=\", \"42-40=\"),\n  @(\"44-10=\", \"80-79=\"),\n  @(\"27+16=\", \"46+5=\"),\n  @(\"26-7=\", \"83-16=\"),\n  @(\"57+22=\", \"33+63=\"),\n  @(\"61+12=\", \"3+63=\"),\n  @(\"25+55=\", \"85-3=\"),\n  @(\"91-72=\", \"65+11=\"),\n  @(\"69-56=\", \"34-28=\"),\n  @(\"73+18=\", \"7+63=\"),\n  @(\"77-23=\", \"46-44=\"),\n  @(\"16+82=\", \"1+75=\"),\n  @(\"55-28=\", \"48-15=\"),\n  @(\"95-63=\", \"96-49=\"),\n  @(\"11-3=\", \"73+24=\"),\n  @(\"39+54=\", \"61-53=\"),\n  @(\"22+11=\", \"82-11=\"),\n  @(\"23+11=\", \"36+17=\"),\n  @(\"75+16=\", \"48-28=\"),\n  @(\"20+70=\", \"54+18=\"),\n  @(\"37+33=\", \"82-41=\"),\n  @(\"74+9=\", \"0+46=\"),\n  @(\"1+88=\", \"97-61=\"),\n  @(\"59+23=\", \"79+6=\"),\n  @(\"4+56=\", \"52-47=\"),\n  @(\"32-12=\", \"45+32=\"),\n  @(\"75-71=\", \"65+28=\"),\n  @(\"83-76=\", \"68-11=\"),\n  @(\"40+42=\", \"56-12=\"),\n  @(\"96-4=\", \"2+77=\"),\n  @(\"90-71=\", \"89-34=\"),\n  @(\"94+2=\", \"46+0=\"),\n  @(\"1+17=\", \"12+30=\"),\n  @(\"24+23=\", \"51-40=\"),\n  @(\"5+7=\", \"95-48=\"),\n  @(\"82-20=\", \"25+29=\"),\n  @(\"60-54=\", \"42-41=\"),\n  @(\"23-16=\", \"37+19=\"),\n  @(\"64-38=\", \"24+24=\"),\n  @(\"9+15=\", \"87-20=\"),\n  @(\"92-75=\", \"83-53=\"),\n  @(\"79+19=\", \"91-54=\"),\n  @(\"67-56=\", \"85-59=\"),\n  @(\"2+66=\", \"36+52=\"),\n  @(\"62-3=\", \"88-65=\"),\n  @(\"30+49=\", \"78-42=\"),\n  @(\"75-61=\", \"59-31=\"),\n  @(\"27+31=\", \"12-11=\"),\n  @(\"56-49=\", \"89-11=\"),\n  @(\"48+6=\", \"74-8=\"),\n  @(\"46-7=\", \"14+57=\"),\n  @(\"52+14=\", \"61-18=\"),\n  @(\"64-57=\", \"63+3=\"),\n  @(\"28+32=\", \"53+18=\"),\n  @(\"70+19=\", \"4+39=\"),\n  @(\"4+93=\", \"68-42=\"),\n  @(\"79-2=\", \"36+5=\"),\n  @(\"32+67=\", \"49+9=\"),\n  @(\"8+58=\", \"37+22=\"),\n  @(\"1+87=\", \"49-48=\"),\n  @(\"58-30=\", \"43+35=\"),\n  @(\"76-5=\", \"39+55=\"),\n  @(\"11+82=\", \"5+51=\"),\n  @(\"78-52=\", \"49-33=\"),\n  @(\"80-73=\", \"10+59=\"),\n  @(\"46+49=\", \"16+45=\"),\n  @(\"52+13=\", \"11+80=\"),\n  @(\"49-44=\", \"79-64=\"),\n  @(\"33-6=\", \"73+12=\"),\n  @(\"20+15=\", \"24-17=\"),\n  @(\"49+38=\", \"39-32=\"),\n  @(\"58-29=\", \"70-13=\"),\n  @(\"19-13=\", \"73-69=\"),\n  @(\"66+3=\", \"0+54=\"),\n  @(\"35+14=\", \"46+44=\"),\n  @(\"41+40=\", \"45-16=\"),\n  @(\"97-41=\", \"81-80=\"),\n)\n\n# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n#   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n#   ReplaceWith, Replace)\n# wdFindContinue = 1, wdReplaceOne = 2 \u2014 every source string below is\n# unique in the document, so a single targeted replace per pair is safe.\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $range = $d.Content\n  $found = $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n  if (-not $found) {\n    throw \"No match found for: $old\"\n  }\n}\n"}
